$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 2023-24 Portland Trail Blazers roster refresh:
#   - Josh Hart, Greg Brown III, Gary Payton II leave the roster
#   - Cam Reddish, Ryan Arcidiacono, Kevin Knox, Matisse Thybulle join
#   - jersey numbers / measurements / exp / college refreshed for the
#     remaining holdovers as rows shift around
# ---------------------------------------------------------------------------

# Final target table (row -> A..K). J (College) and B (No.) can be blank.
$rows = @(
    @{ Row=2;  A=0;  B=17;    C='Shaedon Sharpe';    D='SG'; E='6-6';  F=200; G='May 30, 2003';      H='ca'; I='R';  J='Kentucky';      K='https://www.basketball-reference.com/players/s/sharpsh01.html' },
    @{ Row=3;  A=1;  B=1;     C='Anfernee Simons';   D='SG'; E='6-3';  F=181; G='June 8, 1999';       H='us'; I='4';  J=$null;           K='https://www.basketball-reference.com/players/s/simonan01.html' },
    @{ Row=4;  A=2;  B=9;     C='Jerami Grant';      D='PF'; E='6-8';  F=210; G='March 12, 1994';     H='us'; I='8';  J='Syracuse';      K='https://www.basketball-reference.com/players/g/grantje01.html' },
    @{ Row=5;  A=3;  B=24;    C='Drew Eubanks';      D='C';  E='6-9';  F=245; G='February 1, 1997';   H='us'; I='4';  J='Oregon State';  K='https://www.basketball-reference.com/players/e/eubandr01.html' },
    @{ Row=6;  A=4;  B=27;    C='Jusuf Nurkić';      D='C';  E='6-11'; F=290; G='August 23, 1994';    H='ba'; I='8';  J=$null;           K='https://www.basketball-reference.com/players/n/nurkiju01.html' },
    @{ Row=7;  A=5;  B=0;     C='Damian Lillard';    D='PG'; E='6-2';  F=195; G='July 15, 1990';      H='us'; I='10'; J='Weber State';   K='https://www.basketball-reference.com/players/l/lillada01.html' },
    @{ Row=8;  A=6;  B=2;     C='Trendon Watford';   D='PF'; E='6-9';  F=240; G='November 9, 2000';   H='us'; I='1';  J='LSU';           K='https://www.basketball-reference.com/players/w/watfotr01.html' },
    @{ Row=9;  A=7;  B=34;    C='Jabari Walker';     D='SF'; E='6-9';  F=215; G='July 30, 2002';      H='us'; I='R';  J='Colorado';      K='https://www.basketball-reference.com/players/w/walkeja01.html' },
    @{ Row=10; A=8;  B=10;    C='Nassir Little';     D='SF'; E='6-5';  F=220; G='February 11, 2000';  H='us'; I='3';  J='UNC';           K='https://www.basketball-reference.com/players/l/littlna01.html' },
    @{ Row=11; A=9;  B=26;    C='Justise Winslow';   D='SF'; E='6-6';  F=222; G='March 26, 1996';     H='us'; I='7';  J='Duke';          K='https://www.basketball-reference.com/players/w/winslju01.html' },
    @{ Row=12; A=10; B=6;     C='Keon Johnson';      D='SG'; E='6-5';  F=186; G='March 10, 2002';     H='us'; I='1';  J='Tennessee';     K='https://www.basketball-reference.com/players/j/johnske07.html' },
    @{ Row=13; A=11; B=21;    C='John Butler (TW)';  D='C';  E='7-1';  F=175; G='December 4, 2002';   H='us'; I='R';  J='Florida State'; K='https://www.basketball-reference.com/players/b/butlejo01.html' },
    @{ Row=14; A=12; B=5;     C='Cam Reddish';       D='SF'; E='6-8';  F=218; G='September 1, 1999';  H='us'; I='3';  J='Duke';          K='https://www.basketball-reference.com/players/r/reddica01.html' },
    @{ Row=15; A=13; B=$null; C='Ibou Badji (TW)';   D='C';  E='7-1';  F=240; G='October 13, 2002';   H='sn'; I='R';  J=$null;           K='https://www.basketball-reference.com/players/b/badjiib01.html' },
    @{ Row=16; A=14; B=$null; C='Ryan Arcidiacono';  D='PG'; E='6-3';  F=195; G='March 26, 1994';     H='us'; I='5';  J='Villanova';     K='https://www.basketball-reference.com/players/a/arcidry01.html' },
    @{ Row=17; A=15; B=25;    C='Kevin Knox';        D='SF'; E='6-7';  F=215; G='August 11, 1999';    H='us'; I='4';  J='Kentucky';      K='https://www.basketball-reference.com/players/k/knoxke01.html' },
    @{ Row=18; A=16; B=$null; C='Matisse Thybulle';  D='SG'; E='6-5';  F=201; G='March 4, 1997';      H='us'; I='3';  J='Washington';    K='https://www.basketball-reference.com/players/t/thybuma01.html' }
)

# New row 18 doesn't exist yet: clone formatting (border/font on A, hyperlink
# style on K) from the last existing data row before writing into it.
$ws.Range("A17").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Range("K17").Copy() | Out-Null
$ws.Range("K18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Column I ("Exp") mixes the literal "R" with digit strings ("1".."10") that
# must stay text, not get auto-coerced to numbers - force the whole column
# to Text before writing, then drop back to the Normal style so no stray
# per-cell number format sticks around.
$ws.Range("I2:I18").NumberFormat = "@"

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A

    if ($null -eq $r.B) {
        $ws.Cells.Item($r.Row, 2).ClearContents() | Out-Null
    } else {
        $ws.Cells.Item($r.Row, 2).Value = $r.B
    }

    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
    $ws.Cells.Item($r.Row, 9).Value = $r.I

    if ($null -eq $r.J) {
        $ws.Cells.Item($r.Row, 10).ClearContents() | Out-Null
    } else {
        $ws.Cells.Item($r.Row, 10).Value = $r.J
    }

    $ws.Cells.Item($r.Row, 11).Value = $r.K
}

$ws.Range("I2:I18").Style = "Normal"

# Row 18 (Matisse Thybulle) is brand new - give its bbref URL a live
# hyperlink like every other player row.
$ws.Hyperlinks.Add($ws.Range("K18"), "https://www.basketball-reference.com/players/t/thybuma01.html") | Out-Null
$ws.Range("K17").Copy() | Out-Null
$ws.Range("K18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
